$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.55
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 1.13
$ws.Range("L2").Value = 1.57
$ws.Range("S2").Value = 1.57
$ws.Range("X2").Value = 26
$ws.Range("AF2").Value = 13
$ws.Range("AI2").Value = 34

# Row 3
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 1.14
$ws.Range("L3").Value = 1.62
$ws.Range("N3").Value = 2.88
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 1.62
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 2.38
$ws.Range("S3").Value = 1.53
$ws.Range("X3").Value = 26
$ws.Range("AG3").Value = 12

# Row 4
$ws.Range("H4").Value = 3.8
$ws.Range("M4").Value = 4.4
$ws.Range("Q4").Value = 3.25
$ws.Range("T4").Value = 14.5
$ws.Range("U4").Value = 22
$ws.Range("X4").Value = 25
$ws.Range("AA4").Value = 7.8
$ws.Range("AB4").Value = 12
$ws.Range("AE4").Value = 10.5
$ws.Range("AF4").Value = 11.5
$ws.Range("AH4").Value = 18
$ws.Range("AJ4").Value = 19

# Row 10
$ws.Range("K10").Value = 9
$ws.Range("N10").Value = 2.3
$ws.Range("O10").Value = 1.6

# Row 14
$ws.Range("G14").Value = 2.8
$ws.Range("I14").Value = 2.5
$ws.Range("R14").Value = 1.91
$ws.Range("S14").Value = 1.91

# Row 15
$ws.Range("G15").Value = 2.88
$ws.Range("I15").Value = 2.63
$ws.Range("S15").Value = 1.67

# Row 16
$ws.Range("G16").Value = 2.6
$ws.Range("I16").Value = 2.88
$ws.Range("S16").Value = 1.67

# Row 17
$ws.Range("G17").Value = 1.62
$ws.Range("S17").Value = 1.7

# Row 18
$ws.Range("G18").Value = 2.8
$ws.Range("I18").Value = 2.6

# Row 19
$ws.Range("G19").Value = 2.7
$ws.Range("I19").Value = 2.7

# Row 20
$ws.Range("G20").Value = 2.3

# Row 21
$ws.Range("G21").Value = 2.88
$ws.Range("I21").Value = 2.63
$ws.Range("U21").Value = 13
$ws.Range("AB21").Value = 19

# Row 22
$ws.Range("G22").Value = 2.88
$ws.Range("H22").Value = 2.88
$ws.Range("I22").Value = 2.6

# Row 23
$ws.Range("G23").Value = 2.8
$ws.Range("I23").Value = 2.25

# Row 24
$ws.Range("I24").Value = 2.1

# Row 25
$ws.Range("G25").Value = 3.2
$ws.Range("R25").Value = 2
$ws.Range("S25").Value = 1.73
$ws.Range("AE25").Value = 6.5

# Row 32
$ws.Range("H32").Value = 3.7
$ws.Range("I32").Value = 4.05
$ws.Range("N32").Value = 1.65
$ws.Range("O32").Value = 2
$ws.Range("T32").Value = 7
$ws.Range("U32").Value = 7.5
$ws.Range("V32").Value = 7
$ws.Range("W32").Value = 11.25
$ws.Range("X32").Value = 10.5
$ws.Range("Y32").Value = 17.5
$ws.Range("Z32").Value = 12.5
$ws.Range("AA32").Value = 6.5
$ws.Range("AB32").Value = 11.75
$ws.Range("AC32").Value = 45
$ws.Range("AD32").Value = 250
$ws.Range("AE32").Value = 11.25
$ws.Range("AF32").Value = 19.5
$ws.Range("AG32").Value = 11.25
$ws.Range("AH32").Value = 50
$ws.Range("AI32").Value = 28
$ws.Range("AJ32").Value = 29

# Row 33
$ws.Range("J33").Value = 1.04
$ws.Range("L33").Value = 1.25

# Row 40
$ws.Range("N40").Value = 1.98
$ws.Range("O40").Value = 1.88

Write-Output "Applied all odds updates for 2025-04-29 FlashScore workbook"
